$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column A values (candidate IDs)
$ws.Range("A2").Value = 1111
$ws.Range("A3").Value = 1111
$ws.Range("A4").Value = 1111
$ws.Range("A5").Value = 1112
$ws.Range("A6").Value = 1112
$ws.Range("A7").Value = 1113
$ws.Range("A8").Value = 1113
$ws.Range("A9").Value = 1113
$ws.Range("A10").Value = 1114
$ws.Range("A11").Value = 1114

# Update the selection to A12
$ws.Range("A12").Select()
